$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 335.92856
$ws.Range("I5").Value = 277.9
$ws.Range("J5").Value = 481
$ws.Range("K5").Value = 277.9
$ws.Range("L5").Value = 481
$ws.Range("M5").Value = -162.9
$ws.Range("N5").Value = -711

$ws.Range("H28").Value = 2908.75
$ws.Range("I28").Value = 4783.3335
$ws.Range("K28").Value = 4783.3335
$ws.Range("M28").Value = -4298.3335

$ws.Range("H98").Value = 7517.476
$ws.Range("I98").Value = 8497.777
$ws.Range("J98").Value = 1635.6666
$ws.Range("K98").Value = 8497.777
$ws.Range("L98").Value = 1635.6666
$ws.Range("M98").Value = -6999.777
$ws.Range("N98").Value = -4631.6666

$ws.Range("H112").Value = 1651.5555
$ws.Range("J112").Value = 1651.5555
$ws.Range("L112").Value = 4954.666499999999
$ws.Range("N112").Value = -7170.666499999999

$ws.Range("H122").Value = 7517.476
$ws.Range("I122").Value = 8497.777
$ws.Range("J122").Value = 1635.6666
$ws.Range("K122").Value = 25493.331
$ws.Range("L122").Value = 4906.9998
$ws.Range("M122").Value = -23043.331
$ws.Range("N122").Value = -9806.9998

$ws.Range("H132").Value = 910.9375
$ws.Range("J132").Value = 1438
$ws.Range("L132").Value = 4314
$ws.Range("N132").Value = -9374

$ws.Range("H138").Value = 3055.6052
$ws.Range("J138").Value = 4511.0713
$ws.Range("L138").Value = 13533.2139
$ws.Range("N138").Value = -23813.2139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 37
$ws.Range("I4").Value = 33.75
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 33.75
$ws.Range("L4").Value = 50
$ws.Range("M4").Value = 82.25
$ws.Range("N4").Value = -282

$ws.Range("H32").Value = 3195133.2
$ws.Range("I32").Value = 3299582
$ws.Range("J32").Value = 9444
$ws.Range("K32").Value = 3299582
$ws.Range("L32").Value = 9444
$ws.Range("M32").Value = -3299295
$ws.Range("N32").Value = -10018

$ws.Range("H50").Value = 2654.2
$ws.Range("I50").Value = 2491.6667
$ws.Range("K50").Value = 2491.6667
$ws.Range("M50").Value = -1777.6667

$ws.Range("H61").Value = 6264.1113
$ws.Range("I61").Value = 3943.077
$ws.Range("K61").Value = 3943.077
$ws.Range("M61").Value = -3731.077

$ws.Range("H63").Value = 8860.575999999999
$ws.Range("J63").Value = 10830.27
$ws.Range("L63").Value = 10830.27
$ws.Range("N63").Value = -12202.27

$ws.Range("H66").Value = 8860.575999999999
$ws.Range("J66").Value = 10830.27
$ws.Range("L66").Value = 54151.35000000001
$ws.Range("N66").Value = -61015.35000000001

$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524

$ws.Range("H132").Value = 3489.3547
$ws.Range("I132").Value = 2479.878
$ws.Range("J132").Value = 5460.2383
$ws.Range("K132").Value = 7439.634
$ws.Range("L132").Value = 16380.7149
$ws.Range("M132").Value = -4909.634
$ws.Range("N132").Value = -21440.7149

$ws.Range("H136").Value = 6264.1113
$ws.Range("I136").Value = 3943.077
$ws.Range("K136").Value = 11829.231
$ws.Range("M136").Value = -9279.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 5495
$ws.Range("J15").Value = 4993.3335
$ws.Range("L15").Value = 4993.3335
$ws.Range("N15").Value = -5447.3335

$ws.Range("H86").Value = 334886.34
$ws.Range("J86").Value = 668333.3
$ws.Range("L86").Value = 668333.3
$ws.Range("N86").Value = -670579.3

$ws.Range("H89").Value = 334886.34
$ws.Range("J89").Value = 668333.3
$ws.Range("L89").Value = 3341666.5
$ws.Range("N89").Value = -3352898.5

$ws.Range("H134").Value = 3822.3704
$ws.Range("I134").Value = 2134.4
$ws.Range("J134").Value = 8645.143
$ws.Range("K134").Value = 6403.200000000001
$ws.Range("L134").Value = 25935.429
$ws.Range("M134").Value = -3868.200000000001
$ws.Range("N134").Value = -31005.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1885.7368
$ws.Range("I22").Value = 693.5454999999999
$ws.Range("J22").Value = 3525
$ws.Range("K22").Value = 693.5454999999999
$ws.Range("L22").Value = 3525
$ws.Range("M22").Value = -343.5454999999999
$ws.Range("N22").Value = -4225

$ws.Range("H31").Value = 25645622
$ws.Range("I31").Value = 40003284
$ws.Range("J31").Value = 6938.5713
$ws.Range("K31").Value = 40003284
$ws.Range("L31").Value = 6938.5713
$ws.Range("M31").Value = -40002989
$ws.Range("N31").Value = -7528.5713

$ws.Range("H34").Value = 25645622
$ws.Range("I34").Value = 40003284
$ws.Range("J34").Value = 6938.5713
$ws.Range("K34").Value = 40003284
$ws.Range("L34").Value = 6938.5713
$ws.Range("M34").Value = -40003082
$ws.Range("N34").Value = -7342.5713

$ws.Range("H96").Value = 14489.5
$ws.Range("J96").Value = 14489.5
$ws.Range("L96").Value = 14489.5
$ws.Range("N96").Value = -19981.5

$ws.Range("H107").Value = 835.3333
$ws.Range("I107").Value = 789.75
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 789.75
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 1130.25
$ws.Range("N107").Value = -5040

$ws.Range("H141").Value = 48552.777
$ws.Range("J141").Value = 48552.777
$ws.Range("L141").Value = 48552.777
$ws.Range("N141").Value = -58912.777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 638.61536
$ws.Range("J5").Value = 1059.8
$ws.Range("L5").Value = 3179.4
$ws.Range("N5").Value = -3403.4

$ws.Range("H22").Value = 388.75
$ws.Range("J22").Value = 422
$ws.Range("L22").Value = 1266
$ws.Range("N22").Value = -1604

$ws.Range("H27").Value = 388.75
$ws.Range("J27").Value = 422
$ws.Range("L27").Value = 1266
$ws.Range("N27").Value = -1470

$ws.Range("H32").Value = 808.6667
$ws.Range("I32").Value = 426
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 1278
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -995
$ws.Range("N32").Value = -3566

$ws.Range("H39").Value = 5646.3335
$ws.Range("I39").Value = 5381.143
$ws.Range("J39").Value = 7502.6665
$ws.Range("K39").Value = 16143.429
$ws.Range("L39").Value = 22507.9995
$ws.Range("M39").Value = -15849.429
$ws.Range("N39").Value = -23095.9995

$ws.Range("H57").Value = 1199.5
$ws.Range("I57").Value = 1199.5
$ws.Range("K57").Value = 3598.5
$ws.Range("M57").Value = -3039.5

$ws.Range("H74").Value = 3700
$ws.Range("I74").Value = 3700
$ws.Range("K74").Value = 11100
$ws.Range("M74").Value = -10039

$ws.Range("H77").Value = 3700
$ws.Range("I77").Value = 3700
$ws.Range("K77").Value = 33300
$ws.Range("M77").Value = -27996

$ws.Range("H135").Value = 638.61536
$ws.Range("J135").Value = 1059.8
$ws.Range("L135").Value = 9538.199999999999
$ws.Range("N135").Value = -14608.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 49963.43
$ws.Range("J32").Value = 67748.5
$ws.Range("L32").Value = 67748.5
$ws.Range("N32").Value = -68340.5

$ws.Range("H42").Value = 71645
$ws.Range("J42").Value = 71645
$ws.Range("L42").Value = 71645
$ws.Range("N42").Value = -72615

$ws.Range("H54").Value = 5275
$ws.Range("J54").Value = 5275
$ws.Range("L54").Value = 5275
$ws.Range("N54").Value = -6055

$ws.Range("H70").Value = 8132.6665
$ws.Range("I70").Value = 7949.75
$ws.Range("J70").Value = 8498.5
$ws.Range("K70").Value = 7949.75
$ws.Range("L70").Value = 8498.5
$ws.Range("M70").Value = -7679.75
$ws.Range("N70").Value = -9038.5

$ws.Range("H73").Value = 8132.6665
$ws.Range("I73").Value = 7949.75
$ws.Range("J73").Value = 8498.5
$ws.Range("K73").Value = 7949.75
$ws.Range("L73").Value = 8498.5
$ws.Range("M73").Value = -7013.75
$ws.Range("N73").Value = -10370.5

$ws.Range("H92").Value = 17139.5
$ws.Range("J92").Value = 17139.5
$ws.Range("L92").Value = 17139.5
$ws.Range("N92").Value = -20883.5

$ws.Range("H115").Value = 71645
$ws.Range("J115").Value = 71645
$ws.Range("L115").Value = 71645
$ws.Range("N115").Value = -73995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3164.2683
$ws.Range("I22").Value = 1890.7
$ws.Range("K22").Value = 1890.7
$ws.Range("M22").Value = -1595.7

$ws.Range("H27").Value = 3164.2683
$ws.Range("I27").Value = 1890.7
$ws.Range("K27").Value = 1890.7
$ws.Range("M27").Value = -1783.7

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13363.75
$ws.Range("I41").Value = 16666.334
$ws.Range("J41").Value = 11382.2
$ws.Range("K41").Value = 16666.334
$ws.Range("L41").Value = 11382.2
$ws.Range("M41").Value = -16276.334
$ws.Range("N41").Value = -12162.2

$ws.Range("H92").Value = 48948.75
$ws.Range("J92").Value = 61431.668
$ws.Range("L92").Value = 61431.668
$ws.Range("N92").Value = -66423.66800000001

$ws.Range("H107").Value = 5950
$ws.Range("I107").Value = 5581.8184
$ws.Range("K107").Value = 16745.4552
$ws.Range("M107").Value = -14825.4552

$ws.Range("H126").Value = 85954.5
$ws.Range("I126").Value = 126957
$ws.Range("J126").Value = 3949.5
$ws.Range("K126").Value = 380871
$ws.Range("L126").Value = 11848.5
$ws.Range("M126").Value = -378401
$ws.Range("N126").Value = -16788.5

$ws.Range("H132").Value = 4976.25
$ws.Range("I132").Value = 2721.8
$ws.Range("K132").Value = 8165.400000000001
$ws.Range("M132").Value = -5635.400000000001
